# Apply the changes described by the diff to the GUI naming reference sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: "LabelCheck1" / "label_sample_1-11" -> "Button Stick" / "btn_stick_1-11" ---
$ws.Range("D28").Value = "Button Stick"
$ws.Range("E28").Value = "btn_stick_1-11"

# --- Insert two new rows right after row 29 (the old "Button Sample" row) ---
# This pushes the former rows 30+ down by two (old 30->32, 31->33, 32->34, 33->35, 34->36,
# and the trailing formatted blank row 39 -> 41), matching the target layout.
$ws.Rows.Item(30).Insert()
$ws.Rows.Item(30).Insert()

# Copy the formatting (borders/font) from row 29 onto the two freshly inserted rows
# so they keep style index 2 like their neighbours instead of the generic default style.
$ws.Range("C29:F29").Copy()
$ws.Range("C30:F31").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New row 30 holds the "Button Delete" entry (row 31 remains an empty spacer row).
$ws.Range("D30").Value = "Button Delete"
$ws.Range("E30").Value = "btn_delete_1-11"

# --- New rows describing the camera-related buttons/labels, inserted in the gap
#     between the existing table (now ending at row 36) and the trailing footer row
#     (now shifted to row 41). ---
$ws.Range("D38").Value = "Button check camera "
$ws.Range("E38").Value = "btn_check_cam"

$ws.Range("D39").Value = "Button capture"
$ws.Range("E39").Value = "btn_capture_cam"

$ws.Range("D40").Value = "Button choose camera"
$ws.Range("E40").Value = "btn_choose_cam"

# The footer row (old row 39, now row 41) gains a "Label size" / "lb_size_cam" entry.
$ws.Range("D41").Value = "Label size "
$ws.Range("E41").Value = "lb_size_cam"

# --- Update the view so the active selection / scroll position matches the target ---
# (topLeftCell ends up at A19 in the target; best-effort via the window's scroll
#  position, even though this runtime's writer does not round-trip topLeftCell.)
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
[void]$ws.Range("F39").Select()
